# This workbook's data rows (2-14) get reshuffled: the set of (Fecha, Calidad,
# Volumen, Precio minimo, Precio maximo, Precio promedio ponderado, Unidad de
# comercializacion, Origen, Precio $/Kg, Kg / unidad) tuples stays the same,
# but which tuple sits on which row changes. Columns A, B, C, E, F, G, H, I, J, K
# are identical on every row and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that vary row to row (1-indexed column numbers): D=4 L=12 M=13 N=14
# O=15 P=16 Q=17 R=18 S=19 T=20
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# Snapshot current values for source rows 2..14 before overwriting anything.
# NOTE: use Value2 (not Value) -- in this runtime, reading the bare `.Value`
# property returns a reflection stub instead of the cell's contents.
$snapshot = @{}
for ($r = 2; $r -le 14; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# Maps destination row -> source row (permutation derived from the target diff).
$mapping = @{
    2  = 11
    3  = 14
    4  = 5
    5  = 2
    6  = 3
    7  = 12
    8  = 13
    9  = 8
    10 = 9
    11 = 4
    12 = 6
    13 = 10
    14 = 7
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $srcVals[$c]
    }
}
